# Apply weekly Fruta/Hortaliza (Granada) price-sheet update.
# The diff rewrites the data block in rows 2-14 of the single worksheet;
# each row keeps its static columns (A,B,C,E,F,G,H,I,J) and only the
# variable columns (D,K,L,M,N,O,P,Q,R,S,T) change values/positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value  = 44649            # D2 Fecha
$ws.Cells.Item(2, 11).Value = "Sin especificar" # K2 Variedad
$ws.Cells.Item(2, 12).Value = "Especial"        # L2 Calidad
$ws.Cells.Item(2, 13).Value = 220               # M2 Volumen
$ws.Cells.Item(2, 14).Value = 21600             # N2 Precio minimo
$ws.Cells.Item(2, 15).Value = 21600             # O2 Precio maximo
$ws.Cells.Item(2, 16).Value = 21600             # P2 Precio promedio ponderado
$ws.Cells.Item(2, 17).Value = "$/caja 18 kilos granel" # Q2 Unidad de comercializacion
$ws.Cells.Item(2, 18).Value = "Provincia de Limarí"    # R2 Origen
$ws.Cells.Item(2, 20).Value = 18                # T2 Kg / unidad

# Row 3
$ws.Cells.Item(3, 4).Value  = 44649
$ws.Cells.Item(3, 11).Value = "Sin especificar"
$ws.Cells.Item(3, 12).Value = "Primera"
$ws.Cells.Item(3, 13).Value = 250
$ws.Cells.Item(3, 14).Value = 16200
$ws.Cells.Item(3, 15).Value = 16200
$ws.Cells.Item(3, 16).Value = 16200
$ws.Cells.Item(3, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(3, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(3, 19).Value = 900
$ws.Cells.Item(3, 20).Value = 18

# Row 4
$ws.Cells.Item(4, 4).Value  = 44649
$ws.Cells.Item(4, 11).Value = "Sin especificar"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 180
$ws.Cells.Item(4, 14).Value = 14400
$ws.Cells.Item(4, 15).Value = 14400
$ws.Cells.Item(4, 16).Value = 14400
$ws.Cells.Item(4, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(4, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(4, 19).Value = 800
$ws.Cells.Item(4, 20).Value = 18

# Row 5
$ws.Cells.Item(5, 4).Value  = 44644
$ws.Cells.Item(5, 11).Value = "Sin especificar"
$ws.Cells.Item(5, 12).Value = "Especial"
$ws.Cells.Item(5, 13).Value = 180
$ws.Cells.Item(5, 14).Value = 18000
$ws.Cells.Item(5, 15).Value = 18000
$ws.Cells.Item(5, 16).Value = 18000
$ws.Cells.Item(5, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(5, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 19).Value = 1200

# Row 6
$ws.Cells.Item(6, 4).Value  = 44644
$ws.Cells.Item(6, 11).Value = "Sin especificar"
$ws.Cells.Item(6, 12).Value = "Primera"
$ws.Cells.Item(6, 13).Value = 220
$ws.Cells.Item(6, 14).Value = 13500
$ws.Cells.Item(6, 15).Value = 13500
$ws.Cells.Item(6, 16).Value = 13500
$ws.Cells.Item(6, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(6, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(6, 19).Value = 900

# Row 7
$ws.Cells.Item(7, 4).Value  = 44644
$ws.Cells.Item(7, 12).Value = "Segunda"
$ws.Cells.Item(7, 13).Value = 290
$ws.Cells.Item(7, 14).Value = 12000
$ws.Cells.Item(7, 15).Value = 12000
$ws.Cells.Item(7, 16).Value = 12000
$ws.Cells.Item(7, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(7, 19).Value = 800
$ws.Cells.Item(7, 20).Value = 15

# Row 8
$ws.Cells.Item(8, 4).Value  = 44305
$ws.Cells.Item(8, 11).Value = "Wonderfull"
$ws.Cells.Item(8, 13).Value = 50
$ws.Cells.Item(8, 14).Value = 18000
$ws.Cells.Item(8, 15).Value = 18000
$ws.Cells.Item(8, 16).Value = 18000
$ws.Cells.Item(8, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 1200
$ws.Cells.Item(8, 20).Value = 15

# Row 9
$ws.Cells.Item(9, 4).Value  = 44305
$ws.Cells.Item(9, 11).Value = "Wonderfull"
$ws.Cells.Item(9, 13).Value = 60
$ws.Cells.Item(9, 14).Value = 15000
$ws.Cells.Item(9, 15).Value = 15000
$ws.Cells.Item(9, 16).Value = 15000
$ws.Cells.Item(9, 17).Value = "$/caja 15 kilos granel"
$ws.Cells.Item(9, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(9, 19).Value = 1000
$ws.Cells.Item(9, 20).Value = 15

# Row 10
$ws.Cells.Item(10, 4).Value  = 44285
$ws.Cells.Item(10, 11).Value = "Wonderfull"
$ws.Cells.Item(10, 13).Value = 40
$ws.Cells.Item(10, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(10, 18).Value = "Provincia del Elquí"

# Row 11
$ws.Cells.Item(11, 4).Value  = 44285
$ws.Cells.Item(11, 11).Value = "Wonderfull"
$ws.Cells.Item(11, 13).Value = 90
$ws.Cells.Item(11, 14).Value = 15000
$ws.Cells.Item(11, 15).Value = 15000
$ws.Cells.Item(11, 16).Value = 15000
$ws.Cells.Item(11, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(11, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 19).Value = 1000

# Row 12
$ws.Cells.Item(12, 4).Value  = 44285
$ws.Cells.Item(12, 11).Value = "Wonderfull"
$ws.Cells.Item(12, 13).Value = 75
$ws.Cells.Item(12, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(12, 18).Value = "Provincia del Elquí"

# Row 13
$ws.Cells.Item(13, 4).Value  = 44309
$ws.Cells.Item(13, 13).Value = 40

# Row 14
$ws.Cells.Item(14, 4).Value  = 44309
$ws.Cells.Item(14, 13).Value = 70
